# The author's commit ("Made updates to my PPP, as well as the supporting
# documents.") touched the five person/table shapes on the single
# diagram slide. In the saved OOXML this shows up purely as PowerPoint
# reassigning each touched shape's collaboration bookkeeping id
# (p:nvPr/p:extLst -> p14:modId, schema
# http://schemas.microsoft.com/office/powerpoint/2010/main) -- an
# internal, randomly-generated value PowerPoint stamps on a shape any
# time it is touched/re-saved. It carries no visible/semantic content
# (no position, size, text, or style changes accompany it).
#
# Re-touch each of those five table shapes (by name, so this is robust
# to shape ordering) so each one is registered as modified, without
# altering any visible property.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$touchedNames = @("Table 16", "Table 6", "Table 7", "Table 8", "Table 10")

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($touchedNames -contains $shape.Name) {
        # Re-apply the shape's own name: a harmless write that marks the
        # shape as edited (mirroring whatever minor touch-up the author
        # made) without changing any visible content.
        $shape.Name = $shape.Name
    }
}
